$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.523.81'
$ws.Range('E2').Value = '  +3.11%  '
$ws.Range('D3').Value = '2.437.77'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.536'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '2.435.97'
$ws.Range('E9').Value = '  +1.14%  '
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.351'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.37'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.70%  '
$ws.Range('E15').Value = '  +4.90%  '
$ws.Range('D16').Value = '2.880.91'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').Value = '62.520.59'
$ws.Range('E17').Value = '  +3.19%  '
$ws.Range('D18').Value = '0.0₅0108'
$ws.Range('E18').Value = '  +275.66%  '
$ws.Range('D19').Value = '2.445.26'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.79'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '325.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +9.43%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '631.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.78%  '
$ws.Range('E29').Value = '  +5.28%  '
$ws.Range('D30').Value = '0.0₃0973'
$ws.Range('E30').Value = '  +4.19%  '
$ws.Range('D31').Value = '2.558.50'
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('E33').Value = '  +5.12%  '
$ws.Range('E34').Value = '  +5.73%  '
$ws.Range('E35').Value = '  +2.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.48'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.52%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('E38').Value = '  +3.10%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '152.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.371'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.53'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.55%  '
$ws.Range('E44').Value = '  +4.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('E47').Value = '  +28.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '143.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.40'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.49%  '
$ws.Range('E51').Value = '  +1.48%  '
